$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13: id 23, username "aritz", registration date 28/05/2024.
# The remaining columns (Email, Nombre, Apellidos, Fecha de nacimiento,
# DNI, Acepta politica, Desea recibir info) are left blank, same as the
# other "sparse" sign-up rows already in the sheet (e.g. row 2, 12).
$ws.Cells.Item(13, 1).Value = 23
$ws.Cells.Item(13, 2).Value = "aritz"
$ws.Cells.Item(13, 4).Value = "28/05/2024"
